# Applies the cryptos.xlsx price/volume/coin-name updates described in the diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.035.83"
$ws.Range("E2").Value = "  -1.17%  "

$ws.Range("D3").Value = "3.174.14"
$ws.Range("E3").Value = "  -4.42%  "

$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "590.68"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.13%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "134.51"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -5.80%  "

$ws.Range("E7").Value = "  -0.04%  "

$ws.Range("D8").Value = "3.172.68"
$ws.Range("E8").Value = "  -4.44%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.514"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.91%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.141"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -5.74%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.25"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -5.37%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.454"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.35%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000237"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -4.50%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.97"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.11%  "

$ws.Range("D15").Value = "3.692.64"
$ws.Range("E15").Value = "  -4.55%  "

$ws.Range("E16").Value = "  -1.26%  "

$ws.Range("D17").Value = "3.173.00"
$ws.Range("E17").Value = "  -4.55%  "

$ws.Range("D18").Value = "62.958.93"
$ws.Range("E18").Value = "  -1.43%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.58"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.88%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "461.86"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.76%  "

$ws.Range("E21").Value = "  -1.79%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.698"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -5.13%  "

$ws.Range("E23").Value = "  -4.22%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.48"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.54%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "83.18"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.91%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.12%  "

$ws.Range("E27").Value = "  -0.01%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.67"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.88%  "

$ws.Range("B29").Value = "NEARProtocol"
$ws.Range("C29").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.85"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.78%  "

$ws.Range("B30").Value = "RenderToken"
$ws.Range("C30").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.74"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -6.43%  "

$ws.Range("E31").Value = "  -5.53%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "27.17"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -5.98%  "

$ws.Range("E33").Value = "  -3.99%  "

$ws.Range("E34").Value = "  -6.75%  "

$ws.Range("E35").Value = "  -5.75%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.83"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.92%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "51.44"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.77%  "

$ws.Range("D38").Value = "0.0₃0703"
$ws.Range("E38").Value = "  -6.02%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0388"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.04%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "404.35"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -6.33%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.11"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.78%  "

$ws.Range("E42").Value = "  -4.06%  "

$ws.Range("B43").Value = "dogwifhat"
$ws.Range("C43").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.61"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -5.22%  "

$ws.Range("B44").Value = "Maker"
$ws.Range("C44").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D44").Value = "2.806.28"
$ws.Range("E44").Value = "  -9.73%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.251"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -5.48%  "

$ws.Range("E46").Value = "  +0.02%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.12"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.35%  "

$ws.Range("B48").Value = "Monero"
$ws.Range("C48").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "124.16"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.22%  "

$ws.Range("B49").Value = "InjectiveProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "25.37"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.00%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.111"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.06%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "34.23"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -6.27%  "
